# Slide 24 ("公禱文"), content placeholder shape, second paragraph:
#   "願人都尊袮的名為聖" -> "願人都尊父的名為聖"
# The original single run gets split (as PowerPoint does while the user
# retypes part of the line) into four runs: "願人" | "都" | "尊父的" | "名為聖".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(2, 1)

# Re-apply the paragraph's own Bold formatting onto four character ranges so
# the single run splits into four runs at the boundaries PowerPoint produced:
# "願人" (1-2) | "都" (3) | "尊袮的" (4-6) | "名為聖" (7-9)
$run1 = $para.Characters(1, 2)
$run1.Font.Bold = $run1.Font.Bold

$run2 = $para.Characters(3, 1)
$run2.Font.Bold = $run2.Font.Bold

$run3 = $para.Characters(4, 3)
$run3.Font.Bold = $run3.Font.Bold

$run4 = $para.Characters(7, 3)
$run4.Font.Bold = $run4.Font.Bold

# Now correct the character: 袮 -> 父 (third run, middle character "尊袮的" -> "尊父的")
$run3.Text = "尊父的"
